# "further cleaning to metadata"
#
# - rename the protocol/sample code in column G from E7760 -> E7420
# - give column G (the protocol code column) its own font (Arial 11, black)
#   instead of sharing the generic column style
# - turn column H's literal FALSE boolean cells into a live =FALSE() formula
# - move the selection from H2:H57 down to G2:G57 and scroll the sheet down
#   a bit further

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 57

# --- 1. fix the sample/protocol code shared string (E7760 -> E7420) -------
# Every cell in G2:G57 shares this one string, so rewriting them all collapses
# back onto a single shared-string entry with the corrected text.
$ws.Range("G2:G" + $lastRow).Value = "E7420"

# --- 2. give column G its own distinct font (Arial, 11pt, black) ---------
$gRange = $ws.Range("G2:G" + $lastRow)
$gRange.Font.Name = "Arial"
$gRange.Font.Size = 11
$gRange.Font.Color = 0

# --- 3. replace the literal FALSE values in column H with a =FALSE() -----
#        formula, cell by cell so each keeps its own (non-shared) formula
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Formula = "=FALSE()"
}

# --- 4. move the active selection to G2:G57 and scroll down a bit --------
$ws.Range("G2:G" + $lastRow).Select()
$excel.ActiveWindow.ScrollRow = 28
